# LOM3242.docx restructuring:
# The paragraphs/styles stay in the same sequential order; only the text
# content that occupies each paragraph (and a few runs inside the
# "Avaliação" paragraph) is shuffled around. We capture every "old" value
# first and only then write the "new" values, so the order of writes can
# never clobber a value we still need to read.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Capture all the text currently sitting in the paragraphs that move.
#    (.Range.Text includes the trailing paragraph mark, so trim 1 char.)
# ---------------------------------------------------------------------
function Get-ParaText($index) {
    $r = $d.Paragraphs($index).Range
    return $r.Text.Substring(0, $r.Text.Length - 1)
}

$v6  = Get-ParaText 6    # "A reologia é a ciência..." (PT objetivos)
$v7  = Get-ParaText 7    # "Rheology is the science..." (EN objetivos)
$v9  = Get-ParaText 9    # "5840897 - Clodoaldo Saron"
$v11 = Get-ParaText 11   # "Escoamento de fluidos newtonianos..." (PT resumo)
$v12 = Get-ParaText 12   # "Flow of Newtonian and non-Newtonian..." (EN resumo)
$v14 = Get-ParaText 14   # "1. Introdução. 2. Tensão..." (PT programa)
$v19 = Get-ParaText 19   # bibliography block (SCHRAMM...WHITE, with breaks)

# ---------------------------------------------------------------------
# 2. Capture the three "value" runs inside paragraph 17 (Avaliação),
#    which sit right after the bold labels "Método: ", "Critério: " and
#    "Norma de recuperação: ". Locate them by searching for the (stable,
#    unchanged) bold labels rather than hard-coding character offsets.
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs(17).Range
$endP17 = $p17.End - 1   # exclude the trailing paragraph mark

$fMetodo = $d.Range($p17.Start, $endP17)
[void]$fMetodo.Find.Execute("Método: ")
$metodoEnd = $fMetodo.End

$fCriterio = $d.Range($p17.Start, $endP17)
[void]$fCriterio.Find.Execute("Critério: ")
$criterioStart = $fCriterio.Start
$criterioEnd = $fCriterio.End

$fNorma = $d.Range($p17.Start, $endP17)
[void]$fNorma.Find.Execute("Norma de recuperação: ")
$normaStart = $fNorma.Start
$normaEnd = $fNorma.End

# t1 = old "Método:" value, t2 = old "Critério:" value, t3 = old "Norma de
# recuperação:" value. t1/t2 are followed by a line-break char (Chr 11)
# that must be excluded from the captured text.
$t1 = $d.Range($metodoEnd, $criterioStart - 1).Text
$t2 = $d.Range($criterioEnd, $normaStart - 1).Text
$t3 = $d.Range($normaEnd, $endP17).Text

# ---------------------------------------------------------------------
# 3. Now write the new values. Since everything needed was captured
#    above, the order of these writes does not matter.
# ---------------------------------------------------------------------
$d.Paragraphs(6).Range.Text  = $v11
$d.Paragraphs(7).Range.Text  = $v12
$d.Paragraphs(9).Range.Text  = $v6
$d.Paragraphs(11).Range.Text = $v14
$d.Paragraphs(12).Range.Text = $v7
$d.Paragraphs(14).Range.Text = $t1
$d.Paragraphs(19).Range.Text = $v9

# Re-locate the three value-slots inside paragraph 17 again (their
# character offsets have not shifted, since none of the writes above
# touched paragraph 17), then overwrite them with the next value in the
# chain: Método-slot <- old Critério value, Critério-slot <- old Norma
# value, Norma-slot <- old bibliography block.
$p17b = $d.Paragraphs(17).Range
$endP17b = $p17b.End - 1

$fMetodo2 = $d.Range($p17b.Start, $endP17b)
[void]$fMetodo2.Find.Execute("Método: ")
$metodoEnd2 = $fMetodo2.End

$fCriterio2 = $d.Range($p17b.Start, $endP17b)
[void]$fCriterio2.Find.Execute("Critério: ")
$criterioStart2 = $fCriterio2.Start
$criterioEnd2 = $fCriterio2.End

$fNorma2 = $d.Range($p17b.Start, $endP17b)
[void]$fNorma2.Find.Execute("Norma de recuperação: ")
$normaStart2 = $fNorma2.Start
$normaEnd2 = $fNorma2.End

# Write back-to-front (Norma slot last) so the earlier offsets
# (metodoEnd2 / criterioStart2 / criterioEnd2 / normaStart2) remain valid
# while they are still being used, even though replacing text of a
# different length shifts everything that comes after it.
$d.Range($normaEnd2, $endP17b).Text = $v19
$d.Range($criterioEnd2, $normaStart2 - 1).Text = $t3
$d.Range($metodoEnd2, $criterioStart2 - 1).Text = $t2
